$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Find.Execute("82÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷9=", 1) | Out-Null
$cell = $t.Cell(1, 2)
$cell.Range.Find.Execute("29÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷3=", 1) | Out-Null
$cell = $t.Cell(1, 3)
$cell.Range.Find.Execute("65÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷2=", 1) | Out-Null
$cell = $t.Cell(1, 4)
$cell.Range.Find.Execute("48÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷7=", 1) | Out-Null
$cell = $t.Cell(1, 5)
$cell.Range.Find.Execute("38÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷7=", 1) | Out-Null
$cell = $t.Cell(5, 1)
$cell.Range.Find.Execute("37÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷8=", 1) | Out-Null
$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("46÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷4=", 1) | Out-Null
$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("94÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷4=", 1) | Out-Null
$cell = $t.Cell(5, 4)
$cell.Range.Find.Execute("56÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷6=", 1) | Out-Null
$cell = $t.Cell(5, 5)
$cell.Range.Find.Execute("64÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷6=", 1) | Out-Null
$cell = $t.Cell(9, 1)
$cell.Range.Find.Execute("10÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷5=", 1) | Out-Null
$cell = $t.Cell(9, 2)
$cell.Range.Find.Execute("53÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷9=", 1) | Out-Null
$cell = $t.Cell(9, 3)
$cell.Range.Find.Execute("85÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷3=", 1) | Out-Null
$cell = $t.Cell(9, 4)
$cell.Range.Find.Execute("99÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷8=", 1) | Out-Null
$cell = $t.Cell(9, 5)
$cell.Range.Find.Execute("22÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷7=", 1) | Out-Null
$cell = $t.Cell(13, 1)
$cell.Range.Find.Execute("39÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷5=", 1) | Out-Null
$cell = $t.Cell(13, 2)
$cell.Range.Find.Execute("34÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷9=", 1) | Out-Null
$cell = $t.Cell(13, 3)
$cell.Range.Find.Execute("22÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷9=", 1) | Out-Null
$cell = $t.Cell(13, 4)
$cell.Range.Find.Execute("60÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷4=", 1) | Out-Null
$cell = $t.Cell(13, 5)
$cell.Range.Find.Execute("75÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷9=", 1) | Out-Null
$cell = $t.Cell(17, 1)
$cell.Range.Find.Execute("63÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷6=", 1) | Out-Null
$cell = $t.Cell(17, 2)
$cell.Range.Find.Execute("69÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷5=", 1) | Out-Null
$cell = $t.Cell(17, 3)
$cell.Range.Find.Execute("18÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "96÷3=", 1) | Out-Null
$cell = $t.Cell(17, 4)
$cell.Range.Find.Execute("53÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷4=", 1) | Out-Null
$cell = $t.Cell(17, 5)
$cell.Range.Find.Execute("63÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷2=", 1) | Out-Null

Write-Host "Done applying replacements"
